$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 10001780
$ws.Range("J17").Value = 10001780
$ws.Range("L17").Value = 30005340
$ws.Range("N17").Value = -30005676
$ws.Range("H33").Value = 8336.462
$ws.Range("I33").Value = 8931.5
$ws.Range("K33").Value = 8931.5
$ws.Range("M33").Value = -8702.5
$ws.Range("H92").Value = 125.15385
$ws.Range("J92").Value = 260
$ws.Range("L92").Value = 260
$ws.Range("N92").Value = -2756
$ws.Range("H96").Value = 512.7857
$ws.Range("I96").Value = 572.0909
$ws.Range("J96").Value = 295.33334
$ws.Range("K96").Value = 1716.2727
$ws.Range("L96").Value = 886.0000200000001
$ws.Range("M96").Value = -343.2727
$ws.Range("N96").Value = -3632.00002
$ws.Range("H100").Value = 1323.8
$ws.Range("I100").Value = 1377.3846
$ws.Range("J100").Value = 975.5
$ws.Range("K100").Value = 1377.3846
$ws.Range("L100").Value = 975.5
$ws.Range("M100").Value = -836.3846000000001
$ws.Range("N100").Value = -2057.5
$ws.Range("H111").Value = 6083.4546
$ws.Range("I111").Value = 5466.6665
$ws.Range("J111").Value = 6314.75
$ws.Range("K111").Value = 16399.9995
$ws.Range("L111").Value = 18944.25
$ws.Range("M111").Value = -13332.9995
$ws.Range("N111").Value = -25078.25
$ws.Range("H116").Value = 4428.5713
$ws.Range("H132").Value = 1817.5
$ws.Range("I132").Value = 1705.3334
$ws.Range("K132").Value = 5116.0002
$ws.Range("M132").Value = -2586.0002
$ws.Range("H137").Value = 4312.125
$ws.Range("I137").Value = 2199.4
$ws.Range("K137").Value = 6598.200000000001
$ws.Range("M137").Value = -4048.200000000001
$ws.Range("H138").Value = 4723.0967
$ws.Range("I138").Value = 2302
$ws.Range("J138").Value = 5565.2173
$ws.Range("K138").Value = 6906
$ws.Range("L138").Value = 16695.6519
$ws.Range("M138").Value = -1766
$ws.Range("N138").Value = -26975.6519
$ws.Range("H141").Value = 6072.727
$ws.Range("I141").Value = 5624.706
$ws.Range("K141").Value = 16874.118
$ws.Range("M141").Value = -11694.118

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 104213.85
$ws.Range("I45").Value = 156752.47
$ws.Range("J45").Value = 6642.143
$ws.Range("K45").Value = 156752.47
$ws.Range("L45").Value = 6642.143
$ws.Range("M45").Value = -156375.47
$ws.Range("N45").Value = -7396.143
$ws.Range("H61").Value = 7848.269
$ws.Range("I61").Value = 7793.9585
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 7793.9585
$ws.Range("L61").Value = 8500
$ws.Range("M61").Value = -7581.9585
$ws.Range("N61").Value = -8924
$ws.Range("H105").Value = 31364.666
$ws.Range("J105").Value = 31364.666
$ws.Range("L105").Value = 31364.666
$ws.Range("N105").Value = -38352.666
$ws.Range("H122").Value = 2533.5417
$ws.Range("I122").Value = 2181.238
$ws.Range("K122").Value = 6543.714
$ws.Range("M122").Value = -4093.714
$ws.Range("H132").Value = 7500
$ws.Range("I132").Value = 7500
$ws.Range("K132").Value = 22500
$ws.Range("M132").Value = -19970
$ws.Range("H136").Value = 7848.269
$ws.Range("I136").Value = 7793.9585
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 23381.8755
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -20831.8755
$ws.Range("N136").Value = -30600
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 18500
$ws.Range("J76").Value = 18500
$ws.Range("L76").Value = 18500
$ws.Range("N76").Value = -19130
$ws.Range("H79").Value = 18500
$ws.Range("J79").Value = 18500
$ws.Range("L79").Value = 18500
$ws.Range("N79").Value = -20684
$ws.Range("H107").Value = 1335.9474
$ws.Range("I107").Value = 1348.4615
$ws.Range("K107").Value = 1348.4615
$ws.Range("M107").Value = 571.5385000000001
$ws.Range("H134").Value = 5292
$ws.Range("I134").Value = 4991.7085
$ws.Range("K134").Value = 14975.1255
$ws.Range("M134").Value = -12440.1255
$ws.Range("H138").Value = 59390
$ws.Range("J138").Value = 98780
$ws.Range("L138").Value = 98780
$ws.Range("N138").Value = -109060

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 25728.572
$ws.Range("J50").Value = 25728.572
$ws.Range("L50").Value = 25728.572
$ws.Range("N50").Value = -26978.572
$ws.Range("H62").Value = 9999.333000000001
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 9998
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 9998
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -11246
$ws.Range("H65").Value = 9999.333000000001
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 9998
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 49990
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -56230
$ws.Range("H99").Value = 8614.368
$ws.Range("I99").Value = 7458.9
$ws.Range("K99").Value = 7458.9
$ws.Range("M99").Value = -5960.9
$ws.Range("H126").Value = 8614.368
$ws.Range("I126").Value = 7458.9
$ws.Range("K126").Value = 22376.7
$ws.Range("M126").Value = -19906.7
$ws.Range("H132").Value = 4966.75
$ws.Range("I132").Value = 2689.0667
$ws.Range("K132").Value = 8067.2001
$ws.Range("M132").Value = -5537.2001
$ws.Range("H134").Value = 6198.9287
$ws.Range("I134").Value = 2728.5
$ws.Range("K134").Value = 8185.5
$ws.Range("M134").Value = -5650.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1366.8
$ws.Range("J22").Value = 1451.9565
$ws.Range("L22").Value = 4355.8695
$ws.Range("N22").Value = -4693.8695
$ws.Range("H27").Value = 1366.8
$ws.Range("J27").Value = 1451.9565
$ws.Range("L27").Value = 4355.8695
$ws.Range("N27").Value = -4559.8695
$ws.Range("H81").Value = 1206.0834
$ws.Range("I81").Value = 881.25
$ws.Range("K81").Value = 2643.75
$ws.Range("M81").Value = -1520.75
$ws.Range("H84").Value = 1206.0834
$ws.Range("I84").Value = 881.25
$ws.Range("K84").Value = 7931.25
$ws.Range("M84").Value = -2315.25
$ws.Range("H132").Value = 1547.4
$ws.Range("J132").Value = 1681.3334
$ws.Range("L132").Value = 15132.0006
$ws.Range("N132").Value = -20192.0006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 721.5
$ws.Range("I97").Value = 615.56525
$ws.Range("K97").Value = 615.56525
$ws.Range("M97").Value = -119.56525
$ws.Range("H102").Value = 4741.1665
$ws.Range("I102").Value = 2825.25
$ws.Range("K102").Value = 2825.25
$ws.Range("M102").Value = -1203.25
$ws.Range("H107").Value = 842.1539
$ws.Range("I107").Value = 456.75
$ws.Range("J107").Value = 1013.44446
$ws.Range("K107").Value = 456.75
$ws.Range("L107").Value = 1013.44446
$ws.Range("M107").Value = 1463.25
$ws.Range("N107").Value = -4853.44446
$ws.Range("H126").Value = 3856.9285
$ws.Range("I126").Value = 2666.4443
$ws.Range("J126").Value = 5999.8
$ws.Range("K126").Value = 7999.3329
$ws.Range("L126").Value = 17999.4
$ws.Range("M126").Value = -5529.3329
$ws.Range("N126").Value = -22939.4
$ws.Range("H132").Value = 5352.231
$ws.Range("I132").Value = 4963.8335
$ws.Range("K132").Value = 14891.5005
$ws.Range("M132").Value = -12361.5005
$ws.Range("H133").Value = 142000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 142000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 142000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -152120

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1216
$ws.Range("J82").Value = 827.75
$ws.Range("L82").Value = 827.75
$ws.Range("N82").Value = -1549.75
$ws.Range("H85").Value = 1216
$ws.Range("J85").Value = 827.75
$ws.Range("L85").Value = 827.75
$ws.Range("N85").Value = -3323.75
$ws.Range("H93").Value = 17711.438
$ws.Range("I93").Value = 5241.357
$ws.Range("K93").Value = 5241.357
$ws.Range("M93").Value = -3993.357

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7166.5
$ws.Range("J62").Value = 9998.666999999999
$ws.Range("L62").Value = 9998.666999999999
$ws.Range("N62").Value = -11246.667
$ws.Range("H65").Value = 7166.5
$ws.Range("J65").Value = 9998.666999999999
$ws.Range("L65").Value = 49993.335
$ws.Range("N65").Value = -56233.335
$ws.Range("H100").Value = 1131.3334
$ws.Range("I100").Value = 1000.1111
$ws.Range("K100").Value = 2000.2222
$ws.Range("M100").Value = -1459.2222
$ws.Range("H113").Value = 761.4
$ws.Range("I113").Value = 626.8333
$ws.Range("J113").Value = 963.25
$ws.Range("K113").Value = 1880.4999
$ws.Range("L113").Value = 2889.75
$ws.Range("M113").Value = 289.5001
$ws.Range("N113").Value = -7229.75
$ws.Range("H132").Value = 3888.1755
$ws.Range("I132").Value = 3368.3865
$ws.Range("K132").Value = 10105.1595
$ws.Range("M132").Value = -7575.1595
